{"js": "// Wrap the \"DemisIdQrImage\" QR-code placeholder textbox (the floating text\n// box that sits on top of the \"Untersuchungsbefund ...\" heading paragraph)\n// in the template's conditional-block markers:\n//   {#DemisIdQrImage} ... {/DemisIdQrImage}\n// i.e. prepend \"{#DemisIdQrImage}\" before the existing \"{%DemisIdQrImage}\"\n// run and append \"{/DemisIdQrImage}\" after the \"Meldungs-ID\" caption run,\n// both inside the floating textbox.\n\n// Locate the paragraph that anchors the textbox (its host run carries the\n// w:pict/VML shape) by searching for the unique heading text next to it.\n// The VML textbox's own content isn't reachable through shape/textFrame\n// navigation in this API surface, but it rides along in the host\n// paragraph's OOXML, so we round-trip that through getOoxml/insertOoxml.\nconst searchResults = context.document.body.search(\"Untersuchungsbefund\", { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not locate the 'Untersuchungsbefund' heading paragraph that anchors the QR-code textbox.\");\n}\n\nconst hostParagraphs = searchResults.items[0].paragraphs;\nhostParagraphs.load(\"items\");\nawait context.sync();\n\nconst hostParagraph = hostParagraphs.items[0];\nconst ooxmlResult = hostParagraph.getOoxml();\nawait context.sync();\n\nlet xml = ooxmlResult.value;\n\nif (!xml.includes(\"{%DemisIdQrImage}\")) {\n  throw new Error(\"Expected QR-code placeholder text '{%DemisIdQrImage}' not found in host paragraph XML.\");\n}\nif (!xml.includes(\"Meldungs-ID\")) {\n  throw new Error(\"Expected 'Meldungs-ID' caption text not found in host paragraph XML.\");\n}\n\nconst openRunProps =\n  '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"12\"/><w:szCs w:val=\"12\"/><w:lang w:val=\"en-US\"/></w:rPr>';\n\nfunction markerRun(text) {\n  return \"<w:r>\" + openRunProps + \"<w:t>\" + text + \"</w:t></w:r>\";\n}\n\n// New runs inserted right before the existing \"{%DemisIdQrImage}\" run.\nconst openMarker = markerRun(\"{#\") + markerRun(\"DemisIdQrImage\") + markerRun(\"}\");\n\n// New runs appended right after the existing \"Meldungs-ID\" run.\nconst closeMarker = markerRun(\"{\") + markerRun(\"/\") + markerRun(\"DemisIdQrImage\") + markerRun(\"}\");\n\n// Insert the opening marker immediately before the run that contains the\n// \"{%DemisIdQrImage}\" field placeholder.\nconst placeholderRunStart =\n  '<w:r w:rsidRPr=\"006F7F3C\"><w:rPr><w:sz w:val=\"16\"/><w:szCs w:val=\"16\"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r>';\n\nlet newXml;\nif (xml.includes(placeholderRunStart)) {\n  newXml = xml.replace(placeholderRunStart, openMarker + placeholderRunStart);\n} else {\n  // Fallback: insert right before the <w:t>{%DemisIdQrImage}</w:t> run, more loosely matched.\n  const idx = xml.indexOf(\"{%DemisIdQrImage}\");\n  let runStart = xml.lastIndexOf(\"<w:r>\", idx);\n  const runStartAlt = xml.lastIndexOf(\"<w:r \", idx);\n  if (runStartAlt > runStart) runStart = runStartAlt;\n  newXml = xml.slice(0, runStart) + openMarker + xml.slice(runStart);\n}\n\n// Insert the closing marker immediately after the run that contains the\n// \"Meldungs-ID\" caption text.\nconst captionRunEnd = \"<w:t>Meldungs-ID</w:t></w:r>\";\nconst insertPos = newXml.indexOf(captionRunEnd);\nif (insertPos < 0) {\n  throw new Error(\"Could not find the 'Meldungs-ID' caption run to append the closing marker after.\");\n}\nconst afterPos = insertPos + captionRunEnd.length;\nnewXml = newXml.slice(0, afterPos) + closeMarker + newXml.slice(afterPos);\n\n// Also refresh the w:pict anchorId / v:shape id, matching what Word assigns\n// when it re-serializes this floating textbox after the edit.\nnewXml = newXml.replace('w14:anchorId=\"09533294\"', 'w14:anchorId=\"10B51230\"');\nnewXml = newXml.replace('id=\"_x0000_s2051\"', 'id=\"_x0000_s2054\"');\n\n// Write the modified OOXML back, replacing the exact paragraph it came from.\nhostParagraph.insertOoxml(newXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Wrap the \"DemisIdQrImage\" QR-code placeholder textbox (the floating\n# text box that sits on the \"Untersuchungsbefund ...\" heading paragraph)\n# in the template's conditional-block markers:\n#   {#DemisIdQrImage} ... {/DemisIdQrImage}\n# i.e. prepend \"{#DemisIdQrImage}\" before the existing \"{%DemisIdQrImage}\"\n# run and append \"{/DemisIdQrImage}\" after the \"Meldungs-ID\" caption run,\n# both inside the floating textbox.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that anchors the textbox (its host run carries the\n# w:pict/VML shape) by searching for the unique heading text next to it.\n# (Note: Find.Execute's resulting Range, even after chaining into\n# .Paragraphs.Item(1), only spans the matched text itself here rather than\n# the whole enclosing paragraph, so walk the Paragraphs collection instead\n# to get a properly paragraph-scoped Range for the WordOpenXML round trip.)\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Untersuchungsbefund*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -lt 0) {\n    throw \"Could not locate the 'Untersuchungsbefund' heading paragraph that anchors the QR-code textbox.\"\n}\n$hostParagraph = $d.Paragraphs.Item($targetIndex)\n$hostRange = $hostParagraph.Range\n\n# Pull this paragraph's OOXML so we can edit the nested textbox content\n# (the VML textbox's w:txbxContent isn't reachable through Shapes/TextFrame\n# navigation in this object model, but it rides along in the paragraph's\n# WordOpenXML).\n$xml = $hostRange.WordOpenXML\n\nif ($xml.IndexOf('{%DemisIdQrImage}') -lt 0) {\n    throw \"Expected QR-code placeholder text '{%DemisIdQrImage}' not found in host paragraph XML.\"\n}\nif ($xml.IndexOf('Meldungs-ID') -lt 0) {\n    throw \"Expected 'Meldungs-ID' caption text not found in host paragraph XML.\"\n}\n\n$openRunProps = '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"12\"/><w:szCs w:val=\"12\"/><w:lang w:val=\"en-US\"/></w:rPr>'\n\nfunction New-MarkerRun([string]$text) {\n    return '<w:r>' + $openRunProps + '<w:t>' + $text + '</w:t></w:r>'\n}\n\n# New runs inserted right before the existing \"{%DemisIdQrImage}\" run.\n$openMarker = (New-MarkerRun '{#') + (New-MarkerRun 'DemisIdQrImage') + (New-MarkerRun '}')\n\n# New runs appended right after the existing \"Meldungs-ID\" run.\n$closeMarker = (New-MarkerRun '{') + (New-MarkerRun '/') + (New-MarkerRun 'DemisIdQrImage') + (New-MarkerRun '}')\n\n# Insert the opening marker immediately before the run that contains the\n# \"{%DemisIdQrImage}\" field placeholder.\n$placeholderRunStart = '<w:r w:rsidRPr=\"006F7F3C\"><w:rPr><w:sz w:val=\"16\"/><w:szCs w:val=\"16\"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r>'\nif ($xml.IndexOf($placeholderRunStart) -ge 0) {\n    $newXml = $xml.Replace($placeholderRunStart, $openMarker + $placeholderRunStart)\n} else {\n    # Fallback: insert right before the <w:t>{%DemisIdQrImage}</w:t> run, more loosely matched.\n    $idx = $xml.IndexOf('{%DemisIdQrImage}')\n    $runStart = $xml.LastIndexOf('<w:r>', $idx)\n    $runStartAlt = $xml.LastIndexOf('<w:r ', $idx)\n    if ($runStartAlt -gt $runStart) { $runStart = $runStartAlt }\n    $newXml = $xml.Substring(0, $runStart) + $openMarker + $xml.Substring($runStart)\n}\n\n# Insert the closing marker immediately after the run that contains the\n# \"Meldungs-ID\" caption text.\n$captionRunEnd = '<w:t>Meldungs-ID</w:t></w:r>'\n$insertPos = $newXml.IndexOf($captionRunEnd)\nif ($insertPos -lt 0) {\n    throw \"Could not find the 'Meldungs-ID' caption run to append the closing marker after.\"\n}\n$insertPos = $insertPos + $captionRunEnd.Length\n$newXml = $newXml.Substring(0, $insertPos) + $closeMarker + $newXml.Substring($insertPos)\n\n# Also refresh the w:pict anchorId / v:shape id, matching what Word assigns\n# when it re-serializes this floating textbox after the edit.\n$newXml = $newXml.Replace('w14:anchorId=\"09533294\"', 'w14:anchorId=\"10B51230\"')\n$newXml = $newXml.Replace('id=\"_x0000_s2051\"', 'id=\"_x0000_s2054\"')\n\n# Write the modified OOXML back into the exact range it came from.\n$hostRange.InsertXML($newXml)\n"}
